$d = $word.ActiveDocument

# Locate the Subtitle paragraph ("Assessing international law on working
# prisoners") so the two new Author paragraphs can be inserted right after
# it, before the Date paragraph.
$rng = $d.Content
[void]$rng.Find.Execute("Assessing international law on working prisoners",
                         $true, $false, $false, $false, $false,
                         $true, 1, $false, "", 0)
$subtitlePara = $rng.Paragraphs(1)

# Insert first new paragraph ("Ben Jarman") right after the subtitle.
$subtitlePara.Range.InsertParagraphAfter()
$authorPara1 = $subtitlePara.Next()
$authorPara1.Range.Text = "Ben Jarman"
$authorPara1.Style = "Author"

# Insert second new paragraph ("Catherine Heard") right after the first.
$authorPara1.Range.InsertParagraphAfter()
$authorPara2 = $authorPara1.Next()
$authorPara2.Range.Text = "Catherine Heard"
$authorPara2.Style = "Author"
